# Add "NA" values in column E (duplicate_image_filename) for rows 2-21,
# matching the header D... / E... "duplicate_image_filename" column that
# already exists in A1:K1 (E1 is already "duplicate_image_filename").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
